$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Disambiguate the repeated "10815a" lab blank IDs in column A ---
# Rows 2,4,6,8,10,12 all originally contain "10815a" (tailing/blank entries).
# Rename them to 10815a_1 .. 10815a_6 in row order.
$ws.Range("A2").Value = "10815a_1"
$ws.Range("A4").Value = "10815a_2"
$ws.Range("A6").Value = "10815a_3"
$ws.Range("A8").Value = "10815a_4"
$ws.Range("A10").Value = "10815a_5"
$ws.Range("A12").Value = "10815a_6"

# --- 2. Normalize row shading ---
# Previously every other data row (2,4,6,8,10,12) was shaded (green fill) while
# rows 3,5,7,9,11 were unshaded. After the edit, only the last row (12) keeps the
# shaded look; all the other data rows (2-11) become unshaded, matching the
# formatting already used by row 3.
$xlPasteFormats = -4122

$ws.Range("A3:I3").Copy() | Out-Null
$ws.Range("A2:I2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4:I4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A6:I6").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A8:I8").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A10:I10").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- 3. Widen column A slightly to fit the longer "10815a_N" labels ---
# (target raw width ~9.71 "characters"; ColumnWidth is expressed in the same
# units but gets quantized to the nearest pixel by the engine)
$ws.Columns("A").ColumnWidth = 8.86
